$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 ---
$ws.Range("I2").Value = "SI"
$ws.Range("K2").Value = "2025-12-17 13:31:36"

# --- Row 3 ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 5783861406
$ws.Range("C3").Value = "Oooo Vvvvv"
$ws.Range("D3").Value = 18
$ws.Range("E3").Value = "CERCA 1"
$ws.Range("F3").Value = "pulizie"
$ws.Range("G3").Value = "Veline"
$ws.Range("I3").Value = "NO"
$ws.Range("J3").Value = "2025-12-17 13:31:58"
$ws.Range("L3").Value = 456

# --- Row 4 ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 5783861406
$ws.Range("C4").Value = "Oooo Vvvvv"
$ws.Range("D4").Value = 18
$ws.Range("E4").Value = "CERCA 1"
$ws.Range("F4").Value = "pulizie"
$ws.Range("G4").Value = "LavaWater"
$ws.Range("I4").Value = "NO"
$ws.Range("J4").Value = "2025-12-17 13:32:38"
$ws.Range("L4").Value = 467

# --- Row 5 ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 5783861406
$ws.Range("C5").Value = "Oooo Vvvvv"
$ws.Range("D5").Value = 65
$ws.Range("E5").Value = "Rainusso"
$ws.Range("F5").Value = "appartamento"
$ws.Range("G5").Value = "Mocio con bastone e secchio, Scopa, Paletta, Sacchi del patume per clienti"
$ws.Range("I5").Value = "NO"
$ws.Range("J5").Value = "2025-12-17 13:39:56"
$ws.Range("L5").Value = 474

# --- Row 6 ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 5783861406
$ws.Range("C6").Value = "Oooo Vvvvv"
$ws.Range("D6").Value = 15
$ws.Range("E6").Value = "BUON PASTORE"
$ws.Range("F6").Value = "pulizie"
$ws.Range("G6").Value = "Lavapavimenti, Spugne, LavaWater"
$ws.Range("I6").Value = "NO"
$ws.Range("J6").Value = "2025-12-17 14:24:34"
$ws.Range("L6").Value = 505

# --- Row 7 ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 5783861406
$ws.Range("C7").Value = "Oooo Vvvvv"
$ws.Range("D7").Value = 18
$ws.Range("E7").Value = "CERCA 1"
$ws.Range("F7").Value = "appartamento"
$ws.Range("G7").Value = "Sacchi del patume per clienti, Pastiglie lavastoviglie per clienti"
$ws.Range("I7").Value = "SI"
$ws.Range("J7").Value = "2025-12-17 14:25:18"
$ws.Range("K7").Value = "2025-12-17 15:22:20"
$ws.Range("L7").Value = 522

# --- Row 8 ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 5783861406
$ws.Range("C8").Value = "Oooo Vvvvv"
$ws.Range("D8").Value = 65
$ws.Range("E8").Value = "Rainusso"
$ws.Range("F8").Value = "pulizie"
$ws.Range("G8").Value = "📝 ttttt"
$ws.Range("I8").Value = "NO"
$ws.Range("J8").Value = "2025-12-17 15:11:41"
$ws.Range("L8").Value = 569

# --- Row 9 ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 5783861406
$ws.Range("C9").Value = "Oooo Vvvvv"
$ws.Range("D9").Value = 65
$ws.Range("E9").Value = "Rainusso"
$ws.Range("F9").Value = "pulizie"
$ws.Range("G9").Value = "Lavapavimenti, Spugne"
$ws.Range("I9").Value = "SI"
$ws.Range("J9").Value = "2025-12-17 15:23:05"
$ws.Range("K9").Value = "2025-12-17 15:28:01"
$ws.Range("L9").Value = 591

# --- Row 10 ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 5783861406
$ws.Range("C10").Value = "Oooo Vvvvv"
$ws.Range("D10").Value = 65
$ws.Range("E10").Value = "Rainusso"
$ws.Range("F10").Value = "pulizie"
$ws.Range("G10").Value = "Lavapavimenti, Spugne"
$ws.Range("I10").Value = "SI"
$ws.Range("J10").Value = "2025-12-17 15:29:30"
$ws.Range("K10").Value = "2025-12-17 15:29:33"
$ws.Range("L10").Value = 613

# --- Row 11 ---
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 5783861406
$ws.Range("C11").Value = "Oooo Vvvvv"
$ws.Range("D11").Value = 65
$ws.Range("E11").Value = "Rainusso"
$ws.Range("F11").Value = "pulizie"
$ws.Range("G11").Value = "Spugne, Lavapavimenti"
$ws.Range("I11").Value = "NO"
$ws.Range("J11").Value = "2025-12-17 15:58:00"
$ws.Range("L11").Value = 632
